$wb = $excel.ActiveWorkbook

# Insert the two new sheets right after the first sheet (ClinicalEncounterType).
# Creating "Facility" first and then inserting "Provider" right after it yields
# the final order ClinicalEncounterType, Provider, Facility, ... (and matches
# the sheetId allocation seen in the target workbook: Provider=12, Facility=11).
$firstSheet = $wb.Worksheets.Item(1)

$tmpFacilitySheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $firstSheet)
$tmpFacilitySheet.Name = "Facility"

$tmpProviderSheet = $wb.Worksheets.Add($tmpFacilitySheet)
$tmpProviderSheet.Name = "Provider"

# Re-fetch live references by name: the COM layer's "Add(After:=)" rebinds
# earlier handles to whichever sheet now sits in that slot, so grab fresh
# handles once both sheets exist and are named.
$facilitySheet = $wb.Worksheets.Item("Facility")
$providerSheet = $wb.Worksheets.Item("Provider")

# ================= Facility =================
$facilitySheet.Cells.Item(1,1).Value = "id"
$facilitySheet.Cells.Item(1,2).Value = "uuid"
$facilitySheet.Cells.Item(1,3).Value = "name"
$facilitySheet.Cells.Item(1,4).Value = "is_active"
$facilitySheet.Cells.Item(1,5).Value = "notes"

$facilitySheet.Cells.Item(2,3).Value = "Hospital A"
$facilitySheet.Cells.Item(3,3).Value = "Hospital B"
$facilitySheet.Cells.Item(4,3).Value = "Hospital C"
$facilitySheet.Cells.Item(6,3).Value = "Clinic A"
$facilitySheet.Cells.Item(5,3).Value = "Hospital D"
$facilitySheet.Cells.Item(7,3).Value = "Clinic B"
$facilitySheet.Cells.Item(8,3).Value = "Clinic C"
$facilitySheet.Cells.Item(9,3).Value = "Clinic D"

$facilityNotes = @("note1","note2","note3","note4","note5","note6","note7","note8")
for ($i = 0; $i -lt 8; $i++) {
    $facilitySheet.Cells.Item($i + 2, 5).Value = $facilityNotes[$i]
}

$facilityUuids = @("A1A","B2B","C3C","D4D","E5E","F6F","G7G","H8H")
for ($i = 0; $i -lt 8; $i++) {
    $facilitySheet.Cells.Item($i + 2, 2).Value = $facilityUuids[$i]
}

for ($i = 0; $i -lt 8; $i++) {
    $facilitySheet.Cells.Item($i + 2, 1).Value = $i + 1
    $facilitySheet.Cells.Item($i + 2, 4).Value = $true
}

# ================= Provider =================
$providerSheet.Cells.Item(1,1).Value = "id"
$providerSheet.Cells.Item(1,2).Value = "uuid"
$providerSheet.Cells.Item(1,3).Value = "first_name"
$providerSheet.Cells.Item(1,4).Value = "last_name"
$providerSheet.Cells.Item(1,5).Value = "is_active"
$providerSheet.Cells.Item(1,6).Value = "notes"

$providerSheet.Cells.Item(2,3).Value = "George"
$providerSheet.Cells.Item(2,4).Value = "Washington"
$providerSheet.Cells.Item(3,3).Value = "John"
$providerSheet.Cells.Item(3,4).Value = "Adams"
$providerSheet.Cells.Item(4,3).Value = "Thomas"
$providerSheet.Cells.Item(4,4).Value = "Jefferson"
$providerSheet.Cells.Item(5,3).Value = "James"
$providerSheet.Cells.Item(5,4).Value = "Madison"

$providerUuids = @("A1A","B2B","C3C","D4D")
for ($i = 0; $i -lt 4; $i++) {
    $providerSheet.Cells.Item($i + 2, 2).Value = $providerUuids[$i]
}

$providerNotes = @("note1","note2","note3","note4")
for ($i = 0; $i -lt 4; $i++) {
    $providerSheet.Cells.Item($i + 2, 6).Value = $providerNotes[$i]
}

for ($i = 0; $i -lt 4; $i++) {
    $providerSheet.Cells.Item($i + 2, 1).Value = $i + 1
    $providerSheet.Cells.Item($i + 2, 5).Value = $true
}

# Tab color (yellow), matching the other newly-added "fake data" sheets.
$facilitySheet.Tab.Color = 65535
$providerSheet.Tab.Color = 65535

# Provider is the active/selected tab, matching the source workbook.
$providerSheet.Select()
